$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.106.42"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.959.36"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'246.56"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4892"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("D8").Value = "'44.73"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "'0.2967"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").Value = "'0.06818"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'19.05"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").Value = "'106.59"
$ws.Range("E12").Value = "  -5.07%  "
$ws.Range("D13").Value = "'0.07753"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "1.932.09"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "'5.410"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "'0.7096"
$ws.Range("E16").Value = "  +5.31%  "
$ws.Range("D17").Value = "'283.39"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "30.985.40"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "'0.000007753"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "'13.21"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "2.182.61"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "'5.548"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").Value = "'9.958"
$ws.Range("E26").Value = "  +4.80%  "
$ws.Range("D27").Value = "'168.75"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "'20.00"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "'2.191"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").Value = "'0.1058"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "'1.441"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'4.796"
$ws.Range("E32").Value = "  +17.82%  "
$ws.Range("D33").Value = "'4.515"
$ws.Range("E33").Value = "  +9.06%  "
$ws.Range("D34").Value = "'0.05002"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'0.7666"
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("D36").Value = "'1.166"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("D37").Value = "'0.02049"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "'2.735"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "'2.706"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "'2.135"
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("D41").Value = "'6.425"
$ws.Range("E41").Value = "  +9.34%  "
$ws.Range("D42").Value = "'0.8822"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'109.27"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'73.35"
$ws.Range("E44").Value = "  +5.70%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4450"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "1.017.16"
$ws.Range("E47").Value = "  +20.88%  "
$ws.Range("D48").Value = "'7.465"
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "'9.364"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.2589"
$ws.Range("E51").Value = "  +3.23%  "
